# Update cryptocurrency price/volume figures to the Dec 18 2022 19:xx snapshot.
# Numeric-looking values are written with a leading apostrophe (quote-prefix)
# and then re-styled to "Normal" so they are stored as literal text (matching
# the sheet's existing text-based Price column) without leaving a stray
# number-format style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'249.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'22.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.514"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05638"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'6.463"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.8059"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").Value = "'0.07334"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.03117"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.02914"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.09259"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.001672"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.229"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.04744"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.0005812"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "16OneONEWorstin24h"
$ws.Range("D18").Value = "'0.006423"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.005073"
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").Value = "'3.974"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'3.375"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'2.111"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").Value = "'0.0003000"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Value = "'0.04160"
$ws.Range("D40").Style = "Normal"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.007044"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1042"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003302"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "'0.008679"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005640"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.6803"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.01469"
$ws.Range("D48").Style = "Normal"
